# Update countries & provincias Spain
# Applies the daily COVID-data refresh: updates the "last updated" timestamp,
# re-ranks a handful of countries (their row now carries the freshly
# refreshed numbers while the country that used to occupy the row keeps its
# previous numbers one row down), and refreshes a couple of rows in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 06:51"

# --- Row 15 (Pakistan) : values refreshed in place ---------------------
$ws.Range("B15").Value = 231818
$ws.Range("C15").Value = 3344
$ws.Range("D15").Value = 131649
$ws.Range("E15").Value = 95407
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 50
$ws.Range("H15").Value = 4762

# --- Rows 37-38 : Kazajistan overtakes Ucrania --------------------------
$ws.Range("A37").Value = "Kazajistan"
$ws.Range("B37").Value = 48574
$ws.Range("C37").Value = 1403
$ws.Range("D37").Value = 27334
$ws.Range("E37").Value = 21052
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 188

$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 48500
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 21376
$ws.Range("E38").Value = 25875
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1249

# --- Rows 55-56 : Honduras overtakes Guatemala --------------------------
$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 23943
$ws.Range("C55").Value = 1022
$ws.Range("D55").Value = 2490
$ws.Range("E55").Value = 20814
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 10
$ws.Range("H55").Value = 639

$ws.Range("A56").Value = "Guatemala"
$ws.Range("B56").Value = 23248
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 3382
$ws.Range("E56").Value = 18919
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 947

# --- Row 75 (Australia) : values refreshed in place ---------------------
$ws.Range("E75").Value = 1057
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 106

# --- Rows 78-80 : Kirguistan overtakes Consejo Danes / Senegal ----------
$ws.Range("A78").Value = "Kirguistan"
$ws.Range("B78").Value = 7691
$ws.Range("C78").Value = 314
$ws.Range("D78").Value = 2843
$ws.Range("E78").Value = 4756
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 92

$ws.Range("A79").Value = "Consejo Danes para los Refugiados"
$ws.Range("B79").Value = 7411
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 3184
$ws.Range("E79").Value = 4045
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 182

$ws.Range("A80").Value = "Senegal"
$ws.Range("B80").Value = 7400
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 4870
$ws.Range("E80").Value = 2397
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 133

# --- Rows 186-187 : Butan overtakes Lesoto -------------------------------
$ws.Range("A186").Value = "Butan"
$ws.Range("B186").Value = 80
$ws.Range("C186").Value = 2
$ws.Range("D186").Value = 53
$ws.Range("E186").Value = 27
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

$ws.Range("A187").Value = "Lesoto"
$ws.Range("B187").Value = 79
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 11
$ws.Range("E187").Value = 68
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0
